$wb = $excel.ActiveWorkbook

$wsGof = $wb.Worksheets.Item("gof")
$wsEst = $wb.Worksheets.Item("Estimates 0-1")
$wsMain = $wb.Worksheets.Item("Main effect 0-1")

# gof hunk 0
$wsGof.Range("D2").Value = 24451
$wsGof.Range("F2").Value = 24515
$wsGof.Range("G2").Value = 24696

# gof hunk 1
$wsGof.Range("D3").Value = 24427
$wsGof.Range("F3").Value = 24531
$wsGof.Range("G3").Value = 24825

# Estimates 0-1 hunk 2
$wsEst.Range("B2").Value = -0.035
$wsEst.Range("D2").Value = -0.029
$wsEst.Range("E2").Value = 0.212

# Estimates 0-1 hunk 3
$wsEst.Range("B3").Value = -0.28
$wsEst.Range("D3").Value = -0.231
$wsEst.Range("E3").Value = 27.91

# Estimates 0-1 hunk 4
$wsEst.Range("B4").Value = -0.038
$wsEst.Range("D4").Value = -0.031
$wsEst.Range("E4").Value = 0.331

# Estimates 0-1 hunk 5
$wsEst.Range("B5").Value = -0.109
$wsEst.Range("D5").Value = -0.09
$wsEst.Range("E5").Value = 7.426

# Estimates 0-1 hunk 6
$wsEst.Range("B6").Value = -0.077
$wsEst.Range("C6").Value = 0.086
$wsEst.Range("D6").Value = -0.064
$wsEst.Range("E6").Value = 0.802

# Estimates 0-1 hunk 7
$wsEst.Range("B7").Value = 0.142
$wsEst.Range("D7").Value = 0.117
$wsEst.Range("E7").Value = 1.977

# Estimates 0-1 hunk 8
$wsEst.Range("B8").Value = 0.056
$wsEst.Range("D8").Value = 0.046
$wsEst.Range("E8").Value = 0.314

# Estimates 0-1 hunk 9
$wsEst.Range("B9").Value = 0.037
$wsEst.Range("D9").Value = 0.031
$wsEst.Range("E9").Value = 0.145

# Estimates 0-1 hunk 10
$wsEst.Range("B10").Value = 0.216
$wsEst.Range("D10").Value = 0.178
$wsEst.Range("E10").Value = 5.062

# Estimates 0-1 hunk 11
$wsEst.Range("B11").Value = 0
$wsEst.Range("C11").Value = 0.067
$wsEst.Range("D11").Value = 0
$wsEst.Range("E11").Value = 0

# Estimates 0-1 hunk 12
$wsEst.Range("B12").Value = -0.187
$wsEst.Range("D12").Value = -0.154
$wsEst.Range("E12").Value = 3.875

# Estimates 0-1 hunk 13
$wsEst.Range("B13").Value = -0.062
$wsEst.Range("C13").Value = 0.095
$wsEst.Range("D13").Value = -0.051
$wsEst.Range("E13").Value = 0.426

# Estimates 0-1 hunk 14
$wsEst.Range("B14").Value = -0.21
$wsEst.Range("D14").Value = -0.173
$wsEst.Range("E14").Value = 4.785

# Estimates 0-1 hunk 15
$wsEst.Range("B15").Value = 0.111
$wsEst.Range("C15").Value = 0.099
$wsEst.Range("D15").Value = 0.092
$wsEst.Range("E15").Value = 1.257

# Estimates 0-1 hunk 16
$wsEst.Range("B16").Value = -0.19
$wsEst.Range("D16").Value = -0.157
$wsEst.Range("E16").Value = 3.403

# Estimates 0-1 hunk 17
$wsEst.Range("B17").Value = -0.055
$wsEst.Range("C17").Value = 0.066
$wsEst.Range("D17").Value = -0.045
$wsEst.Range("E17").Value = 0.694

# Estimates 0-1 hunk 18
$wsEst.Range("B18").Value = 0.251
$wsEst.Range("D18").Value = 0.207
$wsEst.Range("E18").Value = 6.981

# Estimates 0-1 hunk 19
$wsEst.Range("B19").Value = 0.285
$wsEst.Range("D19").Value = 0.235
$wsEst.Range("E19").Value = 9

# Estimates 0-1 hunk 20
$wsEst.Range("B20").Value = -0.191
$wsEst.Range("D20").Value = -0.158
$wsEst.Range("E20").Value = 9.804

# Estimates 0-1 hunk 21
$wsEst.Range("B21").Value = 0.385
$wsEst.Range("C21").Value = 0.1
$wsEst.Range("D21").Value = 0.318
$wsEst.Range("E21").Value = 14.823

# Estimates 0-1 hunk 22
$wsEst.Range("B22").Value = 0.047
$wsEst.Range("C22").Value = 0.386
$wsEst.Range("D22").Value = 0.039
$wsEst.Range("E22").Value = 0.015

# Main effect 0-1 hunk 23
$wsMain.Range("B2").Value = -0.384
$wsMain.Range("C2").Value = -0.317

# Main effect 0-1 hunk 24
$wsMain.Range("B3").Value = -0.276
$wsMain.Range("C3").Value = -0.228
